# Implement additional folder column for naming files
#
# Inserts a new "Folder" column between "Url" and "SizeWarc" so the
# existing SizeWarc/SizeLog/Last/State columns shift one place to the
# right (C -> D, D -> E, E -> F, F -> G), then fills in the new column's
# header and the known folder name for the FSO row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at C, pushing SizeWarc (and everything after it)
# one column to the right.
$ws.Columns("C:C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Folder"

# Known folder name for the FSO / www.fso-bremen.de row; the fso-bremen.de
# row (row 3) has no folder value yet, so it is left blank.
$ws.Range("C2").Value = "Forschungsstelle Osteuropa"

# Match the active selection recorded after the edit.
$ws.Range("C2").Select() | Out-Null
